$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.199.74"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.489.72"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'584.17"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'172.57"
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "2.489.42"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "2.941.15"
$ws.Range("D15").Value = "'25.44"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "66.994.54"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "2.486.18"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "'11.02"
$ws.Range("E19").Value = "  -6.30%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  -5.15%  "
$ws.Range("D21").Value = "'349.58"
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("D22").Value = "'4.04"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'68.59"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("E25").Value = "  -4.53%  "
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "2.615.17"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").Value = "'509.03"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "'7.78"
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D33").Value = "'1.23"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "'1.76"
$ws.Range("E34").Value = "  -4.02%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'159.68"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  -7.43%  "
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").Value = "  -4.24%  "
$ws.Range("E40").Value = "  -5.76%  "
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'4.82"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").Value = "'2.37"
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("D46").Value = "'38.70"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("D47").Value = "'142.85"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.514"
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'3.45"
$ws.Range("E49").Value = "  -4.36%  "
$ws.Range("E50").Value = "  -6.40%  "
$ws.Range("D51").Value = "'0.0731"
$ws.Range("E51").Value = "  -0.94%  "
